# "Luận thêm phần đối của cung vô chính diệu"
# Adds 105 new rows (A/B columns, duplicated values) listing the
# "tọa thủ / đồng cung ... tại cung đối Mệnh" combinations for the 14
# main stars (14 singles + C(14,2)=91 pairs = 105 total), appended after
# the existing data in Sheet2 (sheet1.xml), continuing the workbook's
# row-numbering convention (row 4293 is intentionally skipped, mirroring
# the existing gap pattern already present in the sheet, e.g. row 4213).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$stars = @(
    "Tử Vi", "Thiên Cơ", "Thái Dương", "Vũ Khúc", "Thiên Đồng",
    "Liêm Trinh", "Thiên Phủ", "Thái Âm", "Tham Lang", "Cự Môn",
    "Thiên Tướng", "Thiên Lương", "Thất Sát", "Phá Quân"
)

$strings = New-Object System.Collections.ArrayList

foreach ($s in $stars) {
    [void]$strings.Add("$s tọa thủ tại cung đối Mệnh")
}

for ($i = 0; $i -lt $stars.Count; $i++) {
    for ($j = $i + 1; $j -lt $stars.Count; $j++) {
        [void]$strings.Add("$($stars[$i]) đồng cung $($stars[$j]) tại cung đối Mệnh")
    }
}

$startRow = 4294
$endRow = 4398
$rowCount = $endRow - $startRow + 1

$data = New-Object 'object[,]' $rowCount,2
for ($k = 0; $k -lt $rowCount; $k++) {
    $val = $strings[$k]
    $data[$k, 0] = $val
    $data[$k, 1] = $val
}

$targetRange = $ws.Range("A$($startRow):B$($endRow)")
$targetRange.Value = $data

# Restore the tab/viewport state to mirror the author's last on-screen
# position after appending the new rows.
$ws.Range("B4308:B4398").Select()
